$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the sheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

# Column C holds the "Förändrad" (last changed) date for each record.
# Every data row's date serial value needs to move from 45171 (2023-09-02)
# to 45172 (2023-09-03).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
